# Auto-generated edit script: update rolling date window in Quantities sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45571
$ws.Cells.Item(3, 1).Value = 45572
$ws.Cells.Item(4, 1).Value = 45573
$ws.Cells.Item(5, 1).Value = 45574
$ws.Cells.Item(6, 1).Value = 45575
$ws.Cells.Item(7, 1).Value = 45576
$ws.Cells.Item(8, 1).Value = 45577
$ws.Cells.Item(9, 1).Value = 45578
$ws.Cells.Item(10, 1).Value = 45579
$ws.Cells.Item(11, 1).Value = 45580
$ws.Cells.Item(12, 1).Value = 45581
$ws.Cells.Item(13, 1).Value = 45582
$ws.Cells.Item(14, 1).Value = 45583
$ws.Cells.Item(15, 1).Value = 45584
$ws.Cells.Item(16, 1).Value = 45585
$ws.Cells.Item(17, 1).Value = 45586
$ws.Cells.Item(18, 1).Value = 45587
$ws.Cells.Item(19, 1).Value = 45588
$ws.Cells.Item(20, 1).Value = 45589
$ws.Cells.Item(21, 1).Value = 45590
$ws.Cells.Item(22, 1).Value = 45591
$ws.Cells.Item(23, 1).Value = 45592
$ws.Cells.Item(24, 1).Value = 45593
$ws.Cells.Item(25, 1).Value = 45594
$ws.Cells.Item(26, 1).Value = 45595
$ws.Cells.Item(27, 1).Value = 45596
$ws.Cells.Item(28, 1).Value = 45597
$ws.Cells.Item(29, 1).Value = 45598
$ws.Cells.Item(30, 1).Value = 45599
$ws.Cells.Item(31, 1).Value = 45570
$ws.Cells.Item(32, 1).Value = 45569
$ws.Cells.Item(33, 1).Value = 45568
$ws.Cells.Item(34, 1).Value = 45567
$ws.Cells.Item(35, 1).Value = 45566
$ws.Cells.Item(36, 1).Value = 45564
$ws.Cells.Item(36, 3).Value = 0.00170247
$ws.Cells.Item(36, 7).Value = 465.80531254
$ws.Cells.Item(36, 10).Value = 485.38834923
$ws.Cells.Item(37, 1).Value = 45565
$ws.Cells.Item(37, 3).Value = 0.00170247
$ws.Cells.Item(37, 7).Value = 465.80531254
$ws.Cells.Item(37, 10).Value = 485.38834923
$ws.Cells.Item(38, 1).Value = 45563
$ws.Cells.Item(39, 1).Value = 45558
$ws.Cells.Item(39, 3).Value = 0.00004012
$ws.Cells.Item(39, 7).Value = 280.99031254
$ws.Cells.Item(39, 10).Value = 1941.48834923
$ws.Cells.Item(40, 1).Value = 45559
$ws.Cells.Item(40, 3).Value = 0.00004012
$ws.Cells.Item(40, 7).Value = 280.99031254
$ws.Cells.Item(40, 10).Value = 1941.48834923
$ws.Cells.Item(41, 1).Value = 45560
$ws.Cells.Item(41, 2).Value = 116.4121952
$ws.Cells.Item(41, 3).Value = 0.00170247
$ws.Cells.Item(41, 4).Value = 0.00885078
$ws.Cells.Item(41, 5).Value = 0.06933635
$ws.Cells.Item(41, 6).Value = 12792.90181321
$ws.Cells.Item(41, 7).Value = 465.80531254
$ws.Cells.Item(41, 8).Value = 0.24
$ws.Cells.Item(41, 9).Value = 1.7904431
$ws.Cells.Item(41, 10).Value = 485.38834923
$ws.Cells.Item(42, 1).Value = 45561
$ws.Cells.Item(42, 2).Value = 116.4121952
$ws.Cells.Item(42, 3).Value = 0.00170247
$ws.Cells.Item(42, 4).Value = 0.00885078
$ws.Cells.Item(42, 5).Value = 0.06933635
$ws.Cells.Item(42, 6).Value = 12792.90181321
$ws.Cells.Item(42, 7).Value = 465.80531254
$ws.Cells.Item(42, 8).Value = 0.24
$ws.Cells.Item(42, 9).Value = 1.7904431
$ws.Cells.Item(42, 10).Value = 485.38834923
$ws.Cells.Item(43, 1).Value = 45562
$ws.Cells.Item(43, 2).Value = 116.4121952
$ws.Cells.Item(43, 3).Value = 0.00170247
$ws.Cells.Item(43, 4).Value = 0.00885078
$ws.Cells.Item(43, 5).Value = 0.06933635
$ws.Cells.Item(43, 6).Value = 12792.90181321
$ws.Cells.Item(43, 7).Value = 465.80531254
$ws.Cells.Item(43, 8).Value = 0.24
$ws.Cells.Item(43, 9).Value = 1.7904431
$ws.Cells.Item(43, 10).Value = 485.38834923

Write-Host ("Updated cells: " + 81)
Write-Host ("New UsedRange: " + $ws.UsedRange.Address())
